$wb = $excel.ActiveWorkbook

# Update header text on both sheets (they share the same string table entries)
$wsSearch = $wb.Worksheets.Item("ableToSearchProducts")
$wsCart   = $wb.Worksheets.Item("ableToAddToCart")

$wsSearch.Range("A1").Value = "ProductName"
$wsSearch.Range("B1").Value = "ProductPrice"

$wsCart.Range("A1").Value = "ProductName"
$wsCart.Range("B1").Value = "ProductPrice"

# Rename "Qty" header to "Quantity" on the cart sheet
$wsCart.Range("C1").Value = "Quantity"

# Move the saved selection on the search sheet back to A1
$null = $wsSearch.Range("A1").Select()

# Restore the originally active sheet (selecting a range on another sheet
# would otherwise switch the active tab away from it)
$null = $wsCart.Activate()
